# Save report to DB
# Appends a new transaction row (row 4) to Sheet1:
#   A4 = 222 (Ref)
#   B4 = 1500 (Amount)
#   C4 = "7/08/2018" (Date, stored as text like the existing date cells)
# and leaves the selection on B4, matching the author's final cursor
# position after entering the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A4").Value = 222
$ws.Range("B4").Value = 1500

# Entering a date-shaped string via .Value would get auto-converted to a
# serial date number (and pull in a new number-format style). Write it as
# a text formula first, then collapse the formula down to its literal
# value with a values-only paste so the cell ends up as a plain shared
# string cell - exactly like the other "Date" column cells (C2/C3) - with
# no extra style applied.
$ws.Range("C4").Formula = "=""7/08/2018"""
$ws.Range("C4").Copy()
$ws.Range("C4").PasteSpecial(-4163)

$ws.Range("B4").Select()
